# update brochure file link address
# The promotions sheet's brochure links pointed at the old "seapae" folder
# path on GitHub; the repo folder was renamed to "region-1", so both
# brochure hyperlinks need their target (and displayed) URL updated, and
# the second row's link - which previously had plain text instead of a
# working hyperlink - needs to become a real hyperlink too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCoursesFeesUrl = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-1/aibt/AIBT_Courses_Fees_2021_VOL_2.2.pdf"
$newQ4BrochureUrl  = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-1/aibt/AIBTSEAPAE_Q4_Brochure_1OCT-31DEC21_VOL1.1.pdf"

# Row 2 ("AIBT Courses Fees 2021.pdf"): re-point the existing hyperlink at
# the new region-1 address and keep the displayed text (the URL itself) in
# sync with it.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B2").Value = $newCoursesFeesUrl
$ws.Hyperlinks.Add($ws.Range("B2"), $newCoursesFeesUrl)

# Row 3 ("AIBT Region1(SEAPAE) Q4 Promotion.pdf"): this cell only had plain
# text before; turn it into a working hyperlink pointing at the new
# region-1 address as well.
$ws.Range("B3").Value = $newQ4BrochureUrl
$ws.Hyperlinks.Add($ws.Range("B3"), $newQ4BrochureUrl)

# Leave the selection on B3, matching where the edit was made.
$ws.Range("B3").Select() | Out-Null
